$wb = $excel.ActiveWorkbook

# Add the new "country" worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "country"

$data = @(
    @("old", "new"),
    @("Burkina Faso", "Burkina Faso"),
    @("COMORES", "Comores"),
    @("CÔTE D'IVOIRE", "Cote d'Ivoire"),
    @("Equatorial Guinea", "Equatorial Guinea"),
    @("Eswatini", "Eswatini"),
    @("GABON", "Gabon"),
    @("Ghana", "Ghana"),
    @("NIGER", "Niger"),
    @("République du Congo", "Republic of Congo"),
    @("SAO TOME ET PRINCIPE", "Sao Tome & Principe"),
    @("SENEGAL", "Senegal"),
    @("Seychelles", "Seychelles"),
    @("Sierra Leone", "Sierra Leone"),
    @("South Sudan", "South Sudan"),
    @("Uganda", "Uganda")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $data[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $data[$i][1]
}

# Column A was widened to best-fit the longest country name
$newSheet.Columns.Item(1).AutoFit()

$newSheet.Activate()
